$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (row 1, column A)
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 16:33"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1191854
$ws.Range("C4").Value = 3732
$ws.Range("E4").Value = 944481
$ws.Range("G4").Value = 104
$ws.Range("H4").Value = 68702

# Row 19 - Paises Bajos
$ws.Range("F19").Value = 683

# Row 90 - Republica de Yibuti
$ws.Range("B90").Value = 1116
$ws.Range("C90").Value = 4
$ws.Range("D90").Value = 713
$ws.Range("E90").Value = 401

# Row 103 - Sri Lanka
$ws.Range("B103").Value = 721
$ws.Range("C103").Value = 3
$ws.Range("D103").Value = 194
$ws.Range("E103").Value = 519
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 8

# Rows 138/139 swap country order: Sierra Leona <-> Cabo Verde
# Row 138 becomes Cabo Verde
$ws.Range("A138").Value = "Cabo Verde"
$ws.Range("B138").Value = 175
$ws.Range("C138").Value = 10
$ws.Range("D138").Value = 37
$ws.Range("E138").Value = 136
$ws.Range("H138").Value = 2

# Row 139 becomes Sierra Leona
$ws.Range("A139").Value = "Sierra Leona"
$ws.Range("B139").Value = 166
$ws.Range("D139").Value = 29
$ws.Range("E139").Value = 129
$ws.Range("H139").Value = 8

# Row 140 - Liberia
$ws.Range("B140").Value = 166
$ws.Range("C140").Value = 8
$ws.Range("E140").Value = 90

# Row 141 - Birmania
$ws.Range("B141").Value = 161
$ws.Range("C141").Value = 6
$ws.Range("D141").Value = 49

# Rows 198/199 swap country order: San Cristobal y Nieves <-> Burundi
# Row 198 becomes Burundi
$ws.Range("A198").Value = "Burundi"
$ws.Range("D198").Value = 7
$ws.Range("H198").Value = 1

# Row 199 becomes San Cristobal y Nieves
$ws.Range("A199").Value = "San Cristobal y Nieves"
$ws.Range("D199").Value = 8
$ws.Range("H199").Value = 0
